# The sheet "CandidateData_Complete" currently has 9 data rows (rows 1-9,
# row 1 being the header). The edit duplicates the 7 candidate rows
# (rows 3-9) as seven new rows (10-16), carrying over columns A-G only
# (the "Remarks" column H is not copied for the new rows).
#
# We use Copy/Paste (rather than re-typing each value) so that cells whose
# text looks numeric (e.g. phone numbers stored as text such as
# "9766460157") keep their original text storage instead of being
# re-interpreted as numbers when assigned through .Value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3:G9").Copy($ws.Range("A10"))

# The source block A3:G9 is a full rectangle, so a handful of cells that
# are blank in the source (and therefore must stay entirely absent in the
# destination, per the target layout) get copied over as empty cells.
# Clear those specific cells so they don't linger as empty entries.
$ws.Range("C10").ClearContents()
$ws.Range("D10").ClearContents()
$ws.Range("D11").ClearContents()
$ws.Range("F11").ClearContents()
$ws.Range("B15").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("D15").ClearContents()
$ws.Range("F15").ClearContents()
$ws.Range("F16").ClearContents()
